# Refresh the "cryptos" price/volume table (Tue Mar 21 17:14:01 UTC 2023 GitHub Actions run).
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h). Column A (rank index) is unchanged.
# D is written through NumberFormat "@" (Text) first wherever the new value would
# otherwise be auto-parsed as a number by Excel, so the literal text (with any
# trailing zeros, e.g. "6.410") is preserved exactly like the source cell did.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.192.04'
$ws.Range('E2').Value = '  +1.34%  '
$ws.Range('D3').Value = '1.804.54'
$ws.Range('E3').Value = '  +2.07%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '336.75'
$ws.Range('E5').Value = '  -0.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9963'
$ws.Range('E6').Value = '  -0.53%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4522'
$ws.Range('E7').Value = '  +19.56%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3535'
$ws.Range('E8').Value = '  +4.64%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.68'
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.152'
$ws.Range('E10').Value = '  +1.76%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07475'
$ws.Range('E11').Value = '  +2.49%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.89'
$ws.Range('E12').Value = '  -1.62%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.002'
$ws.Range('E13').Value = '  -0.10%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.271'
$ws.Range('E14').Value = '  -0.15%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.273'
$ws.Range('E15').Value = '  -0.26%  '
$ws.Range('D16').Value = '1.801.44'
$ws.Range('E16').Value = '  +1.88%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001084'
$ws.Range('E17').Value = '  +2.14%  '
$ws.Range('E18').Value = '  +0.58%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '82.29'
$ws.Range('E19').Value = '  +1.04%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9977'
$ws.Range('E20').Value = '  -0.53%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.27'
$ws.Range('E21').Value = '  +0.30%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.410'
$ws.Range('E22').Value = '  +0.92%  '
$ws.Range('D23').Value = '28.264.78'
$ws.Range('E23').Value = '  +1.51%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.99'
$ws.Range('E24').Value = '  +1.00%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.397'
$ws.Range('E25').Value = '  +0.37%  '
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '20.72'
$ws.Range('E26').Value = '  +2.38%  '
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.443'
$ws.Range('E27').Value = '  +3.79%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '155.14'
$ws.Range('E28').Value = '  +2.57%  '
$ws.Range('D29').Value = '2.012.10'
$ws.Range('E29').Value = '  +2.17%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.305'
$ws.Range('E30').Value = '  -12.31%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '133.56'
$ws.Range('E31').Value = '  +0.33%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.064'
$ws.Range('E32').Value = '  +0.70%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.959'
$ws.Range('E33').Value = '  +0.40%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.09447'
$ws.Range('E34').Value = '  +7.84%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '12.27'
$ws.Range('E35').Value = '  -1.22%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6801'
$ws.Range('E36').Value = '  +1.14%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02378'
$ws.Range('E37').Value = '  +0.69%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06255'
$ws.Range('E38').Value = '  -0.12%  '
$ws.Range('B39').Value = 'Algorand'
$ws.Range('C39').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2169'
$ws.Range('E39').Value = '  +1.85%  '
$ws.Range('B40').Value = 'InternetComputer(DFINITY)'
$ws.Range('C40').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.207'
$ws.Range('E40').Value = '  +0.31%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.488'
$ws.Range('E41').Value = '  +0.45%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.218'
$ws.Range('E42').Value = '  +0.07%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.195'
$ws.Range('E43').Value = '  +0.80%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.07'
$ws.Range('E44').Value = '  +0.53%  '
$ws.Range('B45').Value = 'Frax'
$ws.Range('C45').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9959'
$ws.Range('E45').Value = '  -0.67%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6158'
$ws.Range('E46').Value = '  +0.32%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.855'
$ws.Range('E47').Value = '  +0.39%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '129.54'
$ws.Range('E48').Value = '  -1.61%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.043'
$ws.Range('E49').Value = '  +0.35%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07119'
$ws.Range('E50').Value = '  -2.39%  '
$ws.Range('B51').Value = 'EOS'
$ws.Range('C51').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.169'
$ws.Range('E51').Value = '  -1.52%  '
